$d = $word.ActiveDocument

# Remove the "Packages and data" Heading 1 section (its bookmark + the
# heading paragraph itself). Everything that was below it shifts up; no
# other paragraph content changes.
$removed = $false
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "packages-and-data") {
        $p = $bm.Range.Paragraphs(1)
        $p.Range.Delete()
        $removed = $true
        break
    }
}

if (-not $removed) {
    # Fallback: locate the heading paragraph by its text if the bookmark
    # could not be found for some reason.
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq "Packages and data`r") {
            $p.Range.Delete()
            break
        }
    }
}
